$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 25 (Log in / login button): Hebrew label changed from "היכנס" to "כניסת משמש"
$ws.Range("C25").Value = "כניסת משמש"

# New row 117: "Store departments" / "מחלקות החנות"
$ws.Range("B117").Value = "Store departments"
$ws.Range("C117").Value = "מחלקות החנות"

# New row 118: "shopping cart" / "עגלת הקניות" (reuses existing Hebrew text also used for "Cart")
$ws.Range("B118").Value = "shopping cart"
$ws.Range("C118").Value = "עגלת הקניות"

# Update selection / view state to match the saved workbook (active cell B118)
$ws.Range("B118").Select()
